# Add test case GetListSyunoSeikyuTest:
# - KaikeiInf: move selection to D2
# - Add two new sheets (SyunoSeikyu, SyunoNyukin) with sample data rows

$wb = $excel.ActiveWorkbook

# --- 1. KaikeiInf selection change ---
$wsKaikeiInf = $wb.Worksheets.Item("KaikeiInf")
$wsKaikeiInf.Range("D2").Select()

# --- 2. Add "SyunoSeikyu" sheet (after the last existing sheet) ---
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSeikyu = $wb.Worksheets.Add($null, $wsLast)
$wsSeikyu.Name = "SyunoSeikyu"

$seikyuHeaders = @(
    "hp_id", "pt_id", "sin_date", "raiin_no", "nyukin_kbn", "seikyu_tensu",
    "seikyu_gaku", "seikyu_detail", "create_date", "create_id", "create_machine",
    "update_date", "update_id", "update_machine", "new_seikyu_tensu",
    "new_seikyu_gaku", "new_seikyu_detail", "adjust_futan", "new_adjust_futan"
)
for ($i = 0; $i -lt $seikyuHeaders.Count; $i++) {
    $wsSeikyu.Cells.Item(1, $i + 1).Value = $seikyuHeaders[$i]
}

$wsSeikyu.Range("A2").Value = 998
$wsSeikyu.Range("B2").Value = 12345
$wsSeikyu.Range("C2").Value = 20180807
$wsSeikyu.Range("D2").Value = 1234321
$wsSeikyu.Range("E2").Value = 1
$wsSeikyu.Range("F2").Value = 1144
$wsSeikyu.Range("G2").Value = 3430
$wsSeikyu.Range("I2").Value = 40413
$wsSeikyu.Range("I2").NumberFormat = "mm:ss.0"
$wsSeikyu.Range("J2").Value = 0
$wsSeikyu.Range("L2").Value = 40428
$wsSeikyu.Range("L2").NumberFormat = "mm:ss.0"
$wsSeikyu.Range("M2").Value = 0
$wsSeikyu.Range("O2").Value = 1144
$wsSeikyu.Range("P2").Value = 3430
$wsSeikyu.Range("R2").Value = 0
$wsSeikyu.Range("S2").Value = 0

$wsSeikyu.Columns.Item(5).ColumnWidth = 10.109375
$wsSeikyu.Columns.Item(6).ColumnWidth = 11.44140625
$wsSeikyu.Columns.Item(7).ColumnWidth = 10.77734375

$wsSeikyu.Range("C2").Select()

# --- 3. Add "SyunoNyukin" sheet (after SyunoSeikyu) ---
$wsNyukin = $wb.Worksheets.Add($null, $wsSeikyu)
$wsNyukin.Name = "SyunoNyukin"

$nyukinHeaders = @(
    "hp_id", "raiin_no", "pt_id", "sin_date", "sort_no", "adjust_futan",
    "nyukin_gaku", "payment_method_cd", "uketuke_sbt", "nyukin_cmt", "is_deleted",
    "create_date", "create_id", "create_machine", "update_date", "update_id",
    "update_machine", "seq_no", "nyukin_date", "nyukinji_tensu", "nyukinji_seikyu",
    "nyukinji_detail"
)
for ($i = 0; $i -lt $nyukinHeaders.Count; $i++) {
    $wsNyukin.Cells.Item(1, $i + 1).Value = $nyukinHeaders[$i]
}

$wsNyukin.Range("A2").Value = 998
$wsNyukin.Range("B2").Value = 1234321
$wsNyukin.Range("C2").Value = 12345
$wsNyukin.Range("D2").Value = 20180807
$wsNyukin.Range("E2").Value = 1
$wsNyukin.Range("F2").Value = 0
$wsNyukin.Range("G2").Value = 0
$wsNyukin.Range("H2").Value = 0
$wsNyukin.Range("I2").Value = 0
$wsNyukin.Range("K2").Value = 0
$wsNyukin.Range("L2").Value = 40330
$wsNyukin.Range("L2").NumberFormat = "mm:ss.0"
$wsNyukin.Range("M2").Value = 0
$wsNyukin.Range("O2").Value = 40330
$wsNyukin.Range("O2").NumberFormat = "mm:ss.0"
$wsNyukin.Range("P2").Value = 0
$wsNyukin.Range("R2").Value = 1
$wsNyukin.Range("S2").Value = 20100601
$wsNyukin.Range("T2").Value = 0
$wsNyukin.Range("U2").Value = 0

$wsNyukin.Range("F6").Select()
